$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("H4").Value  = "Nationstar"
$ws.Range("H5").Value  = "JS(Jiushuo Electronics)"
$ws.Range("H10").Value = "G-Switch (Pinzan)"

$ws.Range("A12").Value = "PCB thickness = 1mm"
$ws.Range("A15").Value = "No.7 large switch or No.9 small switch, choose one of the two, only one needs to be used, the small switch can be SMD soldered, the large switch needs to cut the PCB and soldered manually"
$ws.Range("A17").Value = "Small switch:"
$ws.Range("A42").Value = "Big switch:"
$ws.Range("A76").Value = "When connecting solder joint 9, the working status will be determined by whether the TF card is inserted, and switch SW1 (or SW3) will be inactive."
